$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H1").Value = "test"
